$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for "Relay 1 name" / "Relay 2 name" right before the
# "Reserverd" row (old row 15), shifting everything below down by 2.
$ws.Rows("15:16").Insert()

# New rows: Relay 1 name / Relay 2 name (16 bytes each, carved out of the
# old 96-byte reserved block).
$ws.Range("A15").Value = "Relay 1 name"
$ws.Range("B15").Value = 16
$ws.Range("C15").Formula = "=D14+1"
$ws.Range("D15").Formula = "=B15+C15-1"

$ws.Range("A16").Value = "Relay 2 name"
$ws.Range("B16").Value = 16
$ws.Range("C16").Formula = "=D15+1"
$ws.Range("D16").Formula = "=B16+C16-1"

# "Reserverd" row (now row 17) shrinks from 96 to 64 bytes.
$ws.Range("B17").Value = 64
$ws.Rows("17:17").RowHeight = 14.25

# Re-establish the shared Start/End formulas across the shifted ranges so
# every row below keeps computing contiguously (Language row 10 down to the
# last padding row, now 40).
$ws.Range("C10:C25").Formula = "=D9+1"
$ws.Range("D10:D22").Formula = "=B10+C10-1"
$ws.Range("D23:D38").Formula = "=B23+C23-1"
$ws.Range("D39").Formula = "=B39+C39-1"
$ws.Range("D40").Formula = "=B40+C40-1"
$ws.Range("C26:C38").Formula = "=D25+1"
$ws.Range("C39").Formula = "=D38+1"
$ws.Range("C40").Formula = "=D39+1"

# Stray helper cell that appeared alongside the edit.
$ws.Range("E18").Value = 240

# Selection moved onto the newly inserted block.
$ws.Range("C14").Select()
